$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the previously-empty row 50 with a new time log entry
$ws.Range("A50").Value = 41899
$ws.Range("B50").Value = 0.60347222222222219
$ws.Range("C50").Value = 0.61458333333333337
$ws.Range("D50").Value = 0
$ws.Range("F50").Value = "Coding"

# Update the active selection to reflect where the user ended up editing
$ws.Range("C51").Select()

$wb.Save()
